# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Uva" (La Palmera de La Serena) ahead
# of the existing row 62, pushing the former rows 62-67 down to 64-69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 62-63 (existing 62..67 shift down to 64..69).
$ws.Range("A62:A63").EntireRow.Insert()

# New row 62: Flame Seedless, Provincia del Elquí
$ws.Range("A62").Value = 8
$ws.Range("B62").Value = "Terminal La Palmera de La Serena"
$ws.Range("C62").Value = "Coquimbo"
$ws.Range("D62").Value = 44585
$ws.Range("E62").Value = 4
$ws.Range("F62").Value = "Fruta"
$ws.Range("G62").Value = 100109
$ws.Range("H62").Value = "Uva"
$ws.Range("I62").Value = 100109001
$ws.Range("J62").Value = "Uva"
$ws.Range("K62").Value = "Flame Seedless"
$ws.Range("L62").Value = "Primera"
$ws.Range("M62").Value = 600
$ws.Range("N62").Value = 7500
$ws.Range("O62").Value = 8000
$ws.Range("P62").Value = 7750
$ws.Range("Q62").Value = "`$/caja 15 kilos"
$ws.Range("R62").Value = "Provincia del Elquí"
$ws.Range("S62").Value = 517
$ws.Range("T62").Value = 15

# New row 63: Superior Seedless, Provincia del Elquí
$ws.Range("A63").Value = 8
$ws.Range("B63").Value = "Terminal La Palmera de La Serena"
$ws.Range("C63").Value = "Coquimbo"
$ws.Range("D63").Value = 44585
$ws.Range("E63").Value = 4
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100109
$ws.Range("H63").Value = "Uva"
$ws.Range("I63").Value = 100109001
$ws.Range("J63").Value = "Uva"
$ws.Range("K63").Value = "Superior Seedless"
$ws.Range("L63").Value = "Primera"
$ws.Range("M63").Value = 500
$ws.Range("N63").Value = 10500
$ws.Range("O63").Value = 11000
$ws.Range("P63").Value = 10750
$ws.Range("Q63").Value = "`$/caja 15 kilos"
$ws.Range("R63").Value = "Provincia del Elquí"
$ws.Range("S63").Value = 717
$ws.Range("T63").Value = 15
